$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 189 (pushes existing rows 189-283 down to 191-285)
$ws.Range("A189:A190").EntireRow.Insert()

# New row 189: Pimiento, Zafiro rojo, Primera
$ws.Range("A189").Value = 7
$ws.Range("B189").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C189").Value = "Ñuble"
$ws.Range("D189").Value = 44704
$ws.Range("E189").Value = 16
$ws.Range("F189").Value = 100112002
$ws.Range("G189").Value = "Pimiento"
$ws.Range("H189").Value = "Zafiro rojo"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 30
$ws.Range("K189").Value = 45000
$ws.Range("L189").Value = 45000
$ws.Range("M189").Value = 45000
$ws.Range("N189").Value = "$/caja 15 kilos"
$ws.Range("O189").Value = "Región de Arica y Parinacota"
$ws.Range("P189").Value = 3000
$ws.Range("Q189").Value = 15
$ws.Range("R189").Value = "Hortaliza"

# New row 190: Pimiento, Zafiro verde, Primera
$ws.Range("A190").Value = 7
$ws.Range("B190").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C190").Value = "Ñuble"
$ws.Range("D190").Value = 44704
$ws.Range("E190").Value = 16
$ws.Range("F190").Value = 100112002
$ws.Range("G190").Value = "Pimiento"
$ws.Range("H190").Value = "Zafiro verde"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 50
$ws.Range("K190").Value = 25000
$ws.Range("L190").Value = 25000
$ws.Range("M190").Value = 25000
$ws.Range("N190").Value = "$/caja 15 kilos"
$ws.Range("O190").Value = "Región de Arica y Parinacota"
$ws.Range("P190").Value = 1667
$ws.Range("Q190").Value = 15
$ws.Range("R190").Value = "Hortaliza"
